$p = $ppt.ActivePresentation

function Set-DateField($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "כ""ב/כסלו/תשפ""א"
        }
    }
}

# Update the cached date text on the slide master
Set-DateField $p.SlideMaster

# Update the cached date text on every slide layout
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Set-DateField $layouts.Item($L)
}

# Fix the numbering in the first diagram on slide 1
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(3).TextFrame.TextRange.Characters(2,1).Text = "4"   # Oval 45 (S2 -> S4)
$s1.Shapes.Item(4).TextFrame.TextRange.Characters(2,1).Text = "3"   # Oval 74 (S5 -> S3)
$s1.Shapes.Item(6).TextFrame.TextRange.Characters(2,1).Text = "2"   # Oval 63 (S3 -> S2)
$s1.Shapes.Item(7).TextFrame.TextRange.Characters(2,1).Text = "5"   # Oval 64 (S4 -> S5)

